$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 60
$ws.Range("P2").Value = 2.9
$ws.Range("G4").Value = 1.42
$ws.Range("J4").Value = 5.2
$ws.Range("R4").Value = 1.63
$ws.Range("S4").Value = 2.04
$ws.Range("I7").Value = 2.08
$ws.Range("F8").Value = 2.2
$ws.Range("K8").Value = 60
$ws.Range("H10").Value = 1.38
$ws.Range("P10").Value = 2.04
$ws.Range("F11").Value = 2.06
$ws.Range("K11").Value = 60
$ws.Range("H12").Value = 3.7
$ws.Range("J12").Value = 2.84
$ws.Range("H13").Value = 1.92
$ws.Range("I13").Value = 2.08
$ws.Range("J13").Value = 4.1
$ws.Range("K13").Value = 4.8
$ws.Range("P13").Value = 2.46
$ws.Range("F14").Value = 2.24
$ws.Range("G14").Value = 2.66
$ws.Range("H14").Value = 2.74
$ws.Range("I14").Value = 3.4
$ws.Range("K14").Value = 4.9
$ws.Range("P14").Value = 2.58
$ws.Range("Q14").Value = 1.53
$ws.Range("K15").Value = 60
$ws.Range("F16").Value = 3.15
$ws.Range("G16").Value = 3.3
$ws.Range("H16").Value = 2.5
$ws.Range("J16").Value = 3.2
$ws.Range("Q16").Value = 2.12
$ws.Range("Q19").Value = 2.28
$ws.Range("F20").Value = 1.82
$ws.Range("G20").Value = 1.83
$ws.Range("H20").Value = 5.4
$ws.Range("I20").Value = 5.6
$ws.Range("R20").Value = 1.26
$ws.Range("AD20").Value = 22
$ws.Range("AL20").Value = 50
$ws.Range("O21").Value = 1.45
$ws.Range("F22").Value = 1.84
$ws.Range("J22").Value = 2.68
$ws.Range("H23").Value = 1.28
$ws.Range("P23").Value = 2
$ws.Range("Q23").Value = 1.84
$ws.Range("G24").Value = 2.04
$ws.Range("I24").Value = 5.9
$ws.Range("J24").Value = 3.35
$ws.Range("F26").Value = 3.25
$ws.Range("I26").Value = 2.66
$ws.Range("F29").Value = 2.82
$ws.Range("H29").Value = 2.32
$ws.Range("J29").Value = 2.84
$ws.Range("P29").Value = 1.45
$ws.Range("F30").Value = 1.85
$ws.Range("G30").Value = 2.32
$ws.Range("H30").Value = 3.2
$ws.Range("I30").Value = 4.9
$ws.Range("J30").Value = 3.3
$ws.Range("K30").Value = 7
$ws.Range("P30").Value = 1.78
